# add logger + feature 28/5
# Appends four new address/chat rows (A85:A88) to Sheet1, reusing the same
# "boxed" cell formatting already used by rows 75:84 (style index 5 in the
# original workbook: Arial 10, medium gray border all around, wrap text),
# and moves the viewport/selection the same way the author's Excel session
# ended up (scrolled to A76, active cell D86).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new data rows -------------------------------------------------------
$ws.Range("A85").Value = "Cây số 1 giá bao nhiêu"
$ws.Range("A86").Value = "Dia chj bunkbu xa hoa khah , tp buon ma thuot đaklak , sdt 0369333915 , ten bia
"
$ws.Range("A87").Value = "Tôi muốn lấy thêm 02 bao thì giá bao nhiêu"
$ws.Range("A88").Value = "Chiều dài con dao bao nhiêu vậy bạn"

# Copy the existing "boxed" formatting from the last pre-existing row (A84)
# onto the four freshly added cells so they match the look of rows 75:84.
$ws.Range("A84").Copy()
$ws.Range("A85:A88").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row heights match the wrapped-text autofit Excel would have produced
# (short single-line entries stay at 15.75pt, the long wrapped address goes
# to 39.75pt, same as the existing rows above).
$ws.Rows.Item(85).RowHeight = 15.75
$ws.Rows.Item(86).RowHeight = 39.75
$ws.Rows.Item(87).RowHeight = 15.75
$ws.Rows.Item(88).RowHeight = 15.75

# --- viewport / selection -------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 76
$win.ScrollColumn = 1
$null = $ws.Range("D86").Select()
